# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the newly generated data (gh-pages output update).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 2772
    3  = 737
    4  = 91
    5  = 6662
    6  = 1387
    7  = 16
    9  = 29
    10 = 84
    11 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
